$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content updates -------------------------------------------------

# C2: "Fait mais à modifier pour ajouter la notion de priorité" -> "Fait"
$ws.Range("C2").Value = "Fait"

# B3 text is unchanged ("Minimiser le mouvement dans une même voix ...")
# C3 becomes empty (the stray blank/highlighted cell is removed entirely)
$ws.Range("C3").Clear()

# C5: "en cours" -> "Fait mais à modifier pour ajouter la notion de priorité"
$ws.Range("C5").Value = "Fait mais à modifier pour ajouter la notion de priorité"

# C6 (new): "Fait dans la minimisation"
$ws.Range("C6").Value = "Fait dans la minimisation"

# C7 (new): "Fait dans la minimisation"
$ws.Range("C7").Value = "Fait dans la minimisation"

# C8: " " -> "Fait"
$ws.Range("C8").Value = "Fait"

# --- Formatting -------------------------------------------------------------
# Cells that now read "Fait..." get the same highlighted look already used
# on C2/C10 (solid fill + centered/wrapped alignment). Copy the format from
# C10 (unchanged reference cell) so the theme-based fill colour is preserved.
$ws.Range("C10").Copy() | Out-Null
foreach ($addr in @("C2","C5","C6","C7","C8")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# --- Selection change (C7 is now the active cell, as in the saved file) ----
$ws.Range("C7").Select() | Out-Null
